$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of row 24 down into row 25 (new BOM line item)
$ws.Range("A24:K24").Copy()
$ws.Range("A25:K25").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A25").Formula = "=A24+1"
$ws.Range("B25").Value = "E-Switch"
$ws.Range("C25").Value = "TL3342F160QG/TR"
$ws.Range("D25").Value = "SWITCH TACTILE SPST-NO 0.05A 12V"
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = "Digikey"
$ws.Range("H25").Value = "EG2531CT-ND"
$ws.Range("I25").Value = 0.63
$ws.Range("J25").Formula = "=I25*F25"
